# Fix Training Data Issue (#48)
#
# The "Date" column (BF) held a mangled literal ("6-24-2013-14") that was
# produced by concatenating the game date with the season label. NBA stats
# for games played late at night were attributed to the wrong calendar day,
# so the correct value is the ISO date text "2014-06-24" for every data row
# (rows 2-31; row 1 is the "Date" header and is left untouched).
#
# Assigning the literal string "2014-06-24" straight to Range.Value would be
# auto-interpreted as a date serial by Excel's type inference (since it
# matches a recognized date pattern), which also mints a new number-format
# style on the cell. Neither happened in the real edit (the cell keeps its
# original, unstyled, plain-text shape), so instead the text is produced via
# a formula (forcing a text/string result) in a scratch cell, then copied
# into each date cell with Paste Special - Values, which carries over the
# literal text without re-parsing it as a date and without touching the
# cell's formatting. The scratch cell is cleared afterwards so it leaves no
# trace in the saved sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = "BF"
$correctDate = "2014-06-24"
$scratch = "ZZ1"

$ws.Range($scratch).Formula = '="' + $correctDate + '"'
$ws.Range($scratch).Copy()

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("$col$r").PasteSpecial(-4163)
}

$ws.Range($scratch).ClearContents()
